# Reproduces the "Add files via upload" commit: the single "Countries"
# sheet gets split into 4 tabs (Countries, Sheet1, Sheet2, Sheet3), each
# sheet re-using the same 16 headers / per-country rows via shared
# strings, plus the Countries sheet's own selection/tab-state changes.

$wb = $excel.ActiveWorkbook

$headers = @("Country","R1","R2","R3","R4","R5","R6","R7","R8","R20","R21","R22","R23","R24","R25","R26")

# country name -> [styleName, centerLastRowValues(bool used only for Foremz-style rows), values...]
$atlantis    = @(35,40,15,20,11,7,40,38,40,20,15,40,10,20,10)
$brobdingnag = @(40,45,20,15,15,10,45,29,30,15,15,45,20,10,20)
$carpania    = @(10,20,11,11,11,8,30,11,5,9,7,16,6,3,5)
$dinotopia   = @(10,15,11,7,8,11,9,11,7,5,9,8,8,5,10)
$erewhon     = @(8,11,15,11,6,3,3,3,2,2,1,1,1,1,0)
$foremz      = @(30,10,3,10,9,12,14,16,2,1,0,3,4,5,6)

function Write-Header($ws) {
  for ($c = 1; $c -le 16; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
  }
}

# Writes one country data row. $style is the named cell style ("Good"/
# "Neutral"/"Bad") applied to the row label in column A. When $centerRow
# is $true (the "Foremz" row in the original sheet) the label also gets
# centered and the value cells B:P get centered too (matches cellXfs
# index 5 / 1 from the source workbook).
function Write-Row($ws, $row, $name, $values, $style, $centerRow) {
  $ws.Cells.Item($row, 1).Value = $name
  $ws.Cells.Item($row, 1).Style = $style
  if ($centerRow) {
    $ws.Cells.Item($row, 1).HorizontalAlignment = -4108
  }
  for ($c = 2; $c -le 16; $c++) {
    $ws.Cells.Item($row, $c).Value = $values[$c - 2]
    if ($centerRow) {
      $ws.Cells.Item($row, $c).HorizontalAlignment = -4108
    }
  }
}

# Adds the two trailing blank styled rows (row 6 = "Bad" label only, row
# 7 = "Bad"+centered label with centered B:P) that Sheet1/Sheet2 keep
# from the original sheet's unused rows 6/7.
function Write-BlankTail($ws, $row6, $row7) {
  $ws.Cells.Item($row6, 1).Style = "Bad"
  $ws.Cells.Item($row7, 1).Style = "Bad"
  $ws.Cells.Item($row7, 1).HorizontalAlignment = -4108
  for ($c = 2; $c -le 16; $c++) {
    $ws.Cells.Item($row7, $c).HorizontalAlignment = -4108
  }
}

$countries = $wb.Worksheets.Item("Countries")

# --- Sheet1: Atlantis, Brobdingnag, Carpania, Dinotopia ---
$ws1 = $wb.Worksheets.Add($null, $countries)
$ws1.Name = "Sheet1"
Write-Header $ws1
Write-Row $ws1 2 "Atlantis" $atlantis "Good" $false
Write-Row $ws1 3 "Brobdingnag" $brobdingnag "Good" $false
Write-Row $ws1 4 "Carpania" $carpania "Neutral" $false
Write-Row $ws1 5 "Dinotopia" $dinotopia "Neutral" $false
Write-BlankTail $ws1 6 7
$ws1.Range("A3").Select()

# --- Sheet2: Atlantis, Brobdingnag, Erewhon, Foremz ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"
Write-Header $ws2
Write-Row $ws2 2 "Atlantis" $atlantis "Good" $false
Write-Row $ws2 3 "Brobdingnag" $brobdingnag "Good" $false
Write-Row $ws2 4 "Erewhon" $erewhon "Bad" $false
Write-Row $ws2 5 "Foremz" $foremz "Bad" $true
Write-BlankTail $ws2 6 7
$ws2.Range("G12").Select()

# --- Sheet3: Carpania, Dinotopia, Erewhon, Foremz ---
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Sheet3"
Write-Header $ws3
Write-Row $ws3 2 "Carpania" $carpania "Neutral" $false
Write-Row $ws3 3 "Dinotopia" $dinotopia "Neutral" $false
Write-Row $ws3 4 "Erewhon" $erewhon "Bad" $false
Write-Row $ws3 5 "Foremz" $foremz "Bad" $true
$ws3.Range("D18").Select()

# Countries keeps its data untouched; only the view/selection changes -
# the header-row selection becomes a whole-row range and the tab focus
# moves to the newly-added Sheet3 (last sheet activated below).
$countries.Range("A1:P1").Select()

# Sheet3 is the sheet left active/selected in the saved workbook.
$ws3.Activate()
